$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 533.15
$ws.Range("I28").Value = 548.0526
$ws.Range("J28").Value = 250
$ws.Range("K28").Value = 548.0526
$ws.Range("L28").Value = 250
$ws.Range("M28").Value = -63.05259999999998
$ws.Range("N28").Value = -1220

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H116").Value = 564517.75
$ws.Range("I116").Value = 1253188.1
$ws.Range("J116").Value = 13581.5
$ws.Range("K116").Value = 1253188.1
$ws.Range("L116").Value = 13581.5
$ws.Range("M116").Value = -1249746.1
$ws.Range("N116").Value = -20465.5

$ws.Range("H132").Value = 40163084
$ws.Range("I132").Value = 47811210
$ws.Range("J132").Value = 10444.75
$ws.Range("K132").Value = 143433630
$ws.Range("L132").Value = 31334.25
$ws.Range("M132").Value = -143431100
$ws.Range("N132").Value = -36394.25

$ws.Range("H134").Value = 45128.76
$ws.Range("J134").Value = 45128.76
$ws.Range("L134").Value = 45128.76
$ws.Range("N134").Value = -55268.76

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

$ws.Range("H26").Value = 11626.75
$ws.Range("I26").Value = 3002.3333
$ws.Range("J26").Value = 37500
$ws.Range("K26").Value = 3002.3333
$ws.Range("L26").Value = 37500
$ws.Range("M26").Value = -2672.3333
$ws.Range("N26").Value = -38160

$ws.Range("H32").Value = 6373.0635
$ws.Range("I32").Value = 4358.756
$ws.Range("J32").Value = 10127
$ws.Range("K32").Value = 4358.756
$ws.Range("L32").Value = 10127
$ws.Range("M32").Value = -4071.756
$ws.Range("N32").Value = -10701

$ws.Range("H61").Value = 2038.5
$ws.Range("I61").Value = 1907.7391
$ws.Range("K61").Value = 1907.7391
$ws.Range("M61").Value = -1695.7391

$ws.Range("H74").Value = 2946.0527
$ws.Range("I74").Value = 2271.2856
$ws.Range("K74").Value = 2271.2856
$ws.Range("M74").Value = -1397.2856

$ws.Range("H77").Value = 2946.0527
$ws.Range("I77").Value = 2271.2856
$ws.Range("K77").Value = 11356.428
$ws.Range("M77").Value = -6988.428

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()

$ws.Range("H136").Value = 2038.5
$ws.Range("I136").Value = 1907.7391
$ws.Range("K136").Value = 5723.2173
$ws.Range("M136").Value = -3173.2173

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()

$ws.Range("H123").Value = 30716.875
$ws.Range("J123").Value = 30716.875
$ws.Range("L123").Value = 30716.875
$ws.Range("N123").Value = -40516.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 605.5454999999999
$ws.Range("I22").Value = 276.375
$ws.Range("J22").Value = 1483.3334
$ws.Range("K22").Value = 276.375
$ws.Range("L22").Value = 1483.3334
$ws.Range("M22").Value = 73.625
$ws.Range("N22").Value = -2183.3334

$ws.Range("H35").Value = 19517.428
$ws.Range("I35").Value = 1655.75
$ws.Range("J35").Value = 43333
$ws.Range("K35").Value = 1655.75
$ws.Range("L35").Value = 43333
$ws.Range("M35").Value = -1361.75
$ws.Range("N35").Value = -43921

$ws.Range("H58").Value = 1816.473
$ws.Range("I58").Value = 1556.6418
$ws.Range("K58").Value = 1556.6418
$ws.Range("M58").Value = -1353.6418

$ws.Range("H136").Value = 1816.473
$ws.Range("I136").Value = 1556.6418
$ws.Range("K136").Value = 4669.9254
$ws.Range("M136").Value = -2119.9254

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 592.6070999999999
$ws.Range("I113").Value = 567.5454999999999
$ws.Range("J113").Value = 684.5
$ws.Range("K113").Value = 1702.6365
$ws.Range("L113").Value = 2053.5
$ws.Range("M113").Value = 467.3635000000002
$ws.Range("N113").Value = -6393.5

$ws.Range("H131").Value = 9435055
$ws.Range("J131").Value = 960.3674
$ws.Range("L131").Value = 2881.1022
$ws.Range("N131").Value = -12961.1022

$ws.Range("H132").Value = 2122.3572
$ws.Range("J132").Value = 2801.4443
$ws.Range("L132").Value = 25212.9987
$ws.Range("N132").Value = -30272.9987

$ws.Range("H137").Value = 3025.3684
$ws.Range("J137").Value = 3291.2942
$ws.Range("L137").Value = 9873.882599999999
$ws.Range("N137").Value = -20073.8826

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 102000000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H30").Value = 102000000
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws.Range("H62").Value = 28646.334
$ws.Range("J62").Value = 28646.334
$ws.Range("L62").Value = 28646.334
$ws.Range("N62").Value = -30018.334

$ws.Range("H65").Value = 28646.334
$ws.Range("J65").Value = 28646.334
$ws.Range("L65").Value = 85939.00199999999
$ws.Range("N65").Value = -92803.00199999999

$ws.Range("H113").Value = 1189.6471
$ws.Range("I113").Value = 1171.1
$ws.Range("J113").Value = 1216.1428
$ws.Range("K113").Value = 1171.1
$ws.Range("L113").Value = 1216.1428
$ws.Range("M113").Value = 998.9000000000001
$ws.Range("N113").Value = -5556.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2506.5
$ws.Range("I22").Value = 1467
$ws.Range("J22").Value = 3130.2
$ws.Range("K22").Value = 1467
$ws.Range("L22").Value = 3130.2
$ws.Range("M22").Value = -1172
$ws.Range("N22").Value = -3720.2

$ws.Range("H27").Value = 2506.5
$ws.Range("I27").Value = 1467
$ws.Range("J27").Value = 3130.2
$ws.Range("K27").Value = 1467
$ws.Range("L27").Value = 3130.2
$ws.Range("M27").Value = -1360
$ws.Range("N27").Value = -3344.2

$ws.Range("H63").Value = 40085
$ws.Range("J63").Value = 40085
$ws.Range("L63").Value = 40085
$ws.Range("N63").Value = -41583

$ws.Range("H66").Value = 40085
$ws.Range("J66").Value = 40085
$ws.Range("L66").Value = 120255
$ws.Range("N66").Value = -127743

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 46004.4
$ws.Range("I23").Value = 26666.666
$ws.Range("J23").Value = 75011
$ws.Range("K23").Value = 26666.666
$ws.Range("L23").Value = 75011
$ws.Range("M23").Value = -26437.666
$ws.Range("N23").Value = -75469

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H136").Value = 1348.0322
$ws.Range("I136").Value = 542.55
$ws.Range("J136").Value = 2812.5454
$ws.Range("K136").Value = 1627.65
$ws.Range("L136").Value = 8437.636200000001
$ws.Range("M136").Value = 922.3500000000001
$ws.Range("N136").Value = -13537.6362
